# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) held the literal text "6-15-2013-14" for every
# team row; correct it to the proper ISO date string "2014-06-15".
#
# NumberFormat is forced to Text ("@") before the assignment so the
# ISO-looking string is stored verbatim instead of being auto-converted
# into a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")
$dateRange.NumberFormat = "@"
$dateRange.Value = "2014-06-15"
